$d = $word.ActiveDocument

# 1) Split "A média do semestre será computada com base na relação:M=(P1+2P2)/3"
#    into text + line break + line break + text
$d.Content.Find.Execute(
    "A média do semestre será computada com base na relação:M=(P1+2P2)/3",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A média do semestre será computada com base na relação:^l^lM=(P1+2P2)/3",
    2
)

# 2) Split the recuperação paragraph text into three parts separated by line breaks
$d.Content.Find.Execute(
    "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.A média final, para os alunos em recuperação, será computada com base na relação abaixo:MF=(M+RC)/2",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.^l^lA média final, para os alunos em recuperação, será computada com base na relação abaixo:^l^lMF=(M+RC)/2",
    2
)
